$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.898.19"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "2.543.87"
$ws.Range("E3").Value = "  +3.08%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'567.04"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").Value = "'146.85"
$ws.Range("E6").Value = "  +2.21%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("D9").Value = "2.542.98"
$ws.Range("E9").Value = "  +3.02%  "
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("E11").Value = "  -2.13%  "
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("D14").Value = "'27.30"
$ws.Range("E14").Value = "  +2.98%  "
$ws.Range("D15").Value = "2.999.36"
$ws.Range("E15").Value = "  +3.13%  "
$ws.Range("D16").Value = "62.887.09"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "2.550.75"
$ws.Range("E18").Value = "  +3.39%  "
$ws.Range("D19").Value = "'11.43"
$ws.Range("E19").Value = "  +1.98%  "
$ws.Range("D20").Value = "'336.86"
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").Value = "'6.73"
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "'65.26"
$ws.Range("E24").Value = "  -0.67%  "
$ws.Range("D25").Value = "'1.62"
$ws.Range("E25").Value = "  +8.64%  "
$ws.Range("E26").Value = "  -2.89%  "
$ws.Range("E27").Value = "  +11.34%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("E29").Value = "  +2.90%  "
$ws.Range("D30").Value = "'7.28"
$ws.Range("E30").Value = "  +6.54%  "
$ws.Range("D31").Value = "0.0₃0812"
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("D32").Value = "'1.82"
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("D33").Value = "'178.23"
$ws.Range("E33").Value = "  +0.78%  "
$ws.Range("E34").Value = "  +2.62%  "
$ws.Range("D35").Value = "'406.31"
$ws.Range("E35").Value = "  +10.23%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "'18.95"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D39").Value = "'4.36"
$ws.Range("E39").Value = "  -1.87%  "
$ws.Range("E40").Value = "  +3.76%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").Value = "'39.12"
$ws.Range("E42").Value = "  -3.25%  "
$ws.Range("D43").Value = "'153.00"
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("D45").Value = "'20.75"
$ws.Range("E45").Value = "  +1.00%  "
$ws.Range("E46").Value = "  +1.54%  "
$ws.Range("D47").Value = "'0.0960"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("D48").Value = "'0.0517"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("E49").Value = "  +4.22%  "
$ws.Range("E50").Value = "  +1.22%  "
$ws.Range("E51").Value = "  +0.57%  "
